# Add two new "divmod" rows/solutions (Week_3, Jan 6th & 7th) to the function
# reference sheet — a new table row describing Python's divmod() builtin.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: 함수명 / 동작 / 설명 for divmod()
$ws.Range("A26").Value = "divmod(나누어질 숫자, 나눌 숫자)"
$ws.Range("B26").Value = "나눗셈과 관련된 연산을 한 번에 처리"
$ws.Range("C26").Value = "이 함수는 두 숫자를 입력받아 몫과 나머지를 튜플로 반환합니다."

# Match the row height used throughout the rest of the table.
$ws.Rows.Item(26).RowHeight = 40

# Column A ("함수명") has best-fit width; the new entry is the longest value
# in that column, so the column widens to accommodate it.
$ws.Columns.Item(1).EntireColumn.AutoFit()

# The author scrolled the frozen sheet down so row 12 is the first visible
# row under the frozen header while continuing to work near the bottom.
$excel.ActiveWindow.ScrollRow = 12
